$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "PyCharm / Build / Runtime version / Component" block so that it
#    sits BEFORE the "full_retirement_calculator_app.py" paragraph instead of
#    after it, and drop the stale "By default, all bugs ..." paragraph that
#    used to follow the "Component" heading.
# ---------------------------------------------------------------------------

# Locate the "full_retirement_calculator_app.py" paragraph (the short one,
# not the "Run the program full_retirement_calculator_app.py" bullet later
# in the document) and cut it (this removes the paragraph, mark included).
$appParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "full_retirement_calculator_app.py`r") {
        $appParaIndex = $i
        break
    }
}
$appRange = $d.Paragraphs($appParaIndex).Range
$appRange.Cut()

# Find the "Component" heading paragraph and re-insert the cut paragraph
# immediately after it (i.e. right before whatever now follows "Component").
$componentIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Component`r") {
        $componentIndex = $i
        break
    }
}
$insertPos = $d.Paragraphs($componentIndex + 1).Range.Start
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.Paste()

# Remove the now-orphaned "By default, all bugs ..." paragraph entirely.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("By default, all bugs")) {
        $d.Paragraphs($i).Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) "Priority" value: "Set by management" -> "TBD", highlighted yellow.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Set by management", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TBD", 2) | Out-Null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "TBD`r") {
        $d.Paragraphs($i).Range.HighlightColorIndex = 7
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Drop the stray <w:lastRenderedPageBreak/> marker that used to sit on the
#    "Enter the birth year ..." run (a rendering-only artifact; replacing the
#    run's text with itself regenerates the run without that marker).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Enter the birth year when prompted; press Enter.", $true, $false, `
                         $false, $false, $false, $true, 1, $false, `
                         "Enter the birth year when prompted; press Enter.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "Assigned" value: replace the old placeholder paragraph with
#    "unassigned", highlighted yellow, and drop the empty paragraph that used
#    to trail it at the end of the document.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Current owner of the issue.  Initially, assigned to the Component's lead developer.", `
                         $true, $false, $false, $false, $false, $true, 1, $false, "unassigned", 2) | Out-Null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "unassigned`r") {
        $p = $d.Paragraphs($i)
        $p.Range.HighlightColorIndex = 7
        # Merge away the paragraph's own mark so the trailing empty paragraph
        # that used to follow it disappears, leaving "unassigned" as the
        # final paragraph of the document.
        $markRange = $d.Range($p.Range.End - 1, $p.Range.End)
        $markRange.Delete()
        break
    }
}
